$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the hard-coded AF17:AF20 "prazo de registro" dates with the
# live WORKDAY formula (mirrors the pattern already used by the other
# rows' helper columns), pulling one workday back/forward from AE per
# the $V$2 lead-time parameter.
$ws.Range("AF17").Formula = '=IFERROR(WORKDAY(AE17,($V$2-1)),"-")'
$ws.Range("AF18").Formula = '=IFERROR(WORKDAY(AE18,($V$2-1)),"-")'
$ws.Range("AF19").Formula = '=IFERROR(WORKDAY(AE19,($V$2-1)),"-")'
$ws.Range("AF20").Formula = '=IFERROR(WORKDAY(AE20,($V$2-1)),"-")'

# Scroll/selection moved from N1/AJ26 to S1/AJ22.
$win = $excel.ActiveWindow
$win.ScrollColumn = 19
$win.ScrollRow = 1
$ws.Range("AJ22").Select() | Out-Null
